# The source data for "Hortaliza, Mercado Mayorista Lo Valledor de Santiago -
# Zapallo italiano" is a daily price log, sorted most-recent-first. A new
# day's record was added at the top of the data block (row 276, right after
# the header row), pushing every existing record down by one row (the former
# last record, row 365, becomes row 366). The sheet dimension grows from
# A1:R365 to A1:R366 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 276:365 down to 277:366 by inserting a blank row at 276.
$ws.Rows('276:276').Insert()

# Populate the newly inserted row with the new daily record.
$ws.Range('A276').Value = 6
$ws.Range('B276').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C276').Value = 'Metropolitana'
$ws.Range('D276').Value = 44642
$ws.Range('E276').Value = 13
$ws.Range('F276').Value = 100112032
$ws.Range('G276').Value = 'Zapallo italiano'
$ws.Range('H276').Value = 'Sin especificar'
$ws.Range('I276').Value = 'Primera'
$ws.Range('J276').Value = 450
$ws.Range('K276').Value = 11000
$ws.Range('L276').Value = 12000
$ws.Range('M276').Value = 11356
$ws.Range('N276').Value = '$/caja 60 unidades'
$ws.Range('O276').Value = 'Limache'
$ws.Range('P276').Value = 189
$ws.Range('Q276').Value = 60
$ws.Range('R276').Value = 'Hortaliza'
